# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table switches from the deck's custom "Table_0" style to the
#    built-in "{9C24CEB3-F66E-48A8-84FB-5BF74C57BBB5}" table style.
# 2) The presentation theme (ppt/theme/theme2.xml, the theme actually used by
#    the slide master / all slides) and the notes-master theme
#    (ppt/theme/theme1.xml) had their colour schemes swapped: the deck's
#    visible theme goes from the "Red Violet" / Integral palette back to the
#    stock Office palette. We reproduce that by rewriting the 12 theme colour
#    slots (accessible on any slide's ThemeColorScheme, which is backed by
#    the shared slide-master theme part) to the stock Office RGB values.

function RGBVal([int]$r, [int]$g, [int]$b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{9C24CEB3-F66E-48A8-84FB-5BF74C57BBB5}", $true)

# --- 2. Swap the theme colour scheme back to the stock Office palette ----
# Order of ThemeColorScheme.Item(n): dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    (RGBVal 0x00 0x00 0x00),  # dk1
    (RGBVal 0xFF 0xFF 0xFF),  # lt1
    (RGBVal 0x44 0x54 0x6A),  # dk2
    (RGBVal 0xE7 0xE6 0xE6),  # lt2
    (RGBVal 0x5B 0x9B 0xD5),  # accent1
    (RGBVal 0xED 0x7D 0x31),  # accent2
    (RGBVal 0xA5 0xA5 0xA5),  # accent3
    (RGBVal 0xFF 0xC0 0x00),  # accent4
    (RGBVal 0x44 0x72 0xC4),  # accent5
    (RGBVal 0x70 0xAD 0x47),  # accent6
    (RGBVal 0x05 0x63 0xC1),  # hlink
    (RGBVal 0x95 0x4F 0x72)   # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
